$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): id, abbreviation, shortDisplayName, displayName, name,
#     nickname, location, standingSummary, created ---
# B1 already carries the bold "header" style (s=1) from the old "Created"
# header; copy that formatting across the rest of the header row so the new
# header cells C1:I1 pick up the same style (avoids creating new font/xf
# entries in styles.xml).
$ws.Range("B1").Copy()
$ws.Range("C1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 1).Value = "id"
$ws.Cells.Item(1, 2).Value = "abbreviation"
$ws.Cells.Item(1, 3).Value = "shortDisplayName"
$ws.Cells.Item(1, 4).Value = "displayName"
$ws.Cells.Item(1, 5).Value = "name"
$ws.Cells.Item(1, 6).Value = "nickname"
$ws.Cells.Item(1, 7).Value = "location"
$ws.Cells.Item(1, 8).Value = "standingSummary"
$ws.Cells.Item(1, 9).Value = "created"

# --- Old sample data is no longer correct; drop it ---
$ws.Rows.Item(3).Delete()
$ws.Range("A2").ClearContents()

# --- New row 2: only the "created" column is filled in ---
$ws.Cells.Item(2, 9).Value = "Fri Jul  5 06:51:04 2024"

# --- Column widths: the two old text columns go back to default width,
#     column I (created) gets a custom width ---
$ws.Columns.Item(1).ColumnWidth = 8.6171875
$ws.Columns.Item(2).ColumnWidth = 8.6171875
$ws.Columns.Item(9).ColumnWidth = 21.1

# --- Selection back to the top-left cell ---
$ws.Range("A1").Select()
